# Scheduled-runner update: refresh currentAveragePrice / Leve price-profit
# columns (H:N) on the ALC, ARM, BSM, CRP, CUL and LTW leve-profit sheets
# with newly pulled market data. Only numeric cell values change (a couple
# of rows gain/lose their profit columns when prices move to/from zero);
# no rows, columns, formulas or styles are added or removed otherwise.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("H132").Value = 2531.3076
$ws.Range("I132").Value = 1173.4546
$ws.Range("K132").Value = 3520.3638
$ws.Range("M132").Value = -990.3638000000001
$ws.Range("H138").Value = 12860.921
$ws.Range("J138").Value = 13249.543
$ws.Range("L138").Value = 39748.629
$ws.Range("N138").Value = -50028.629

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 100000
$ws.Range("J135").Value = 100000
$ws.Range("L135").Value = 100000
$ws.Range("N135").Value = -110140
$ws.Range("H139").Value = 100000
$ws.Range("J139").Value = 100000
$ws.Range("L139").Value = 100000
$ws.Range("N139").Value = -110280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3669.6667
$ws.Range("I99").Value = 3504.5
$ws.Range("K99").Value = 3504.5
$ws.Range("M99").Value = -2006.5
$ws.Range("H107").Value = 1800
$ws.Range("I107").Value = 1720
$ws.Range("K107").Value = 1720
$ws.Range("M107").Value = 200
$ws.Range("H135").Value = 63027
$ws.Range("J135").Value = 63027
$ws.Range("L135").Value = 63027
$ws.Range("N135").Value = -73167

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2767.0833
$ws.Range("I31").Value = 1525.75
$ws.Range("J31").Value = 5249.75
$ws.Range("K31").Value = 1525.75
$ws.Range("L31").Value = 5249.75
$ws.Range("M31").Value = -1230.75
$ws.Range("N31").Value = -5839.75
$ws.Range("H34").Value = 2767.0833
$ws.Range("I34").Value = 1525.75
$ws.Range("J34").Value = 5249.75
$ws.Range("K34").Value = 1525.75
$ws.Range("L34").Value = 5249.75
$ws.Range("M34").Value = -1323.75
$ws.Range("N34").Value = -5653.75
$ws.Range("H50").Value = 55992
$ws.Range("J50").Value = 55992
$ws.Range("L50").Value = 55992
$ws.Range("N50").Value = -57242
$ws.Range("H51").Value = 52540.75
$ws.Range("J51").Value = 52540.75
$ws.Range("L51").Value = 52540.75
$ws.Range("N51").Value = -54012.75
$ws.Range("H59").Value = 34257.75
$ws.Range("J59").Value = 63563.5
$ws.Range("L59").Value = 63563.5
$ws.Range("N59").Value = -65853.5
$ws.Range("H60").Value = 26663.334
$ws.Range("I60").Value = 14997.5
$ws.Range("K60").Value = 14997.5
$ws.Range("M60").Value = -14486.5
$ws.Range("H61").Value = 52540.75
$ws.Range("J61").Value = 52540.75
$ws.Range("L61").Value = 52540.75
$ws.Range("N61").Value = -53236.75
$ws.Range("H68").Value = 79295
$ws.Range("J68").Value = 79295
$ws.Range("L68").Value = 79295
$ws.Range("N68").Value = -80793
$ws.Range("H71").Value = 79295
$ws.Range("J71").Value = 79295
$ws.Range("L71").Value = 237885
$ws.Range("N71").Value = -245373
$ws.Range("H74").Value = 59983.25
$ws.Range("J74").Value = 59983.25
$ws.Range("L74").Value = 59983.25
$ws.Range("N74").Value = -61731.25
$ws.Range("H77").Value = 59983.25
$ws.Range("J77").Value = 59983.25
$ws.Range("L77").Value = 179949.75
$ws.Range("N77").Value = -188685.75
$ws.Range("H122").Value = 1385.2
$ws.Range("I122").Value = 1385.2
$ws.Range("K122").Value = 4155.6
$ws.Range("M122").Value = -1705.6
$ws.Range("H132").Value = 203779.4
$ws.Range("I132").Value = 335666
$ws.Range("J132").Value = 5949.5
$ws.Range("K132").Value = 1006998
$ws.Range("L132").Value = 17848.5
$ws.Range("M132").Value = -1004468
$ws.Range("N132").Value = -22908.5
$ws.Range("H134").Value = 500
$ws.Range("I134").Value = 500
$ws.Range("K134").Value = 1500
$ws.Range("M134").Value = 1035

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 71.5625
$ws.Range("J12").Value = 70.8
$ws.Range("L12").Value = 212.4
$ws.Range("N12").Value = -558.4
$ws.Range("H34").Value = 3129.111
$ws.Range("J34").Value = 5230
$ws.Range("L34").Value = 15690
$ws.Range("N34").Value = -15858
$ws.Range("H39").Value = 15500
$ws.Range("J39").Value = 15500
$ws.Range("L39").Value = 46500
$ws.Range("N39").Value = -47088
$ws.Range("H55").Value = 6966.6665
$ws.Range("J55").Value = 9250
$ws.Range("L55").Value = 27750
$ws.Range("N55").Value = -28104
$ws.Range("H80").Value = 13600
$ws.Range("I80").Value = 8000
$ws.Range("J80").Value = 15000
$ws.Range("K80").Value = 24000
$ws.Range("L80").Value = 45000
$ws.Range("M80").Value = -23064
$ws.Range("N80").Value = -46872
$ws.Range("H83").Value = 13600
$ws.Range("I83").Value = 8000
$ws.Range("J83").Value = 15000
$ws.Range("K83").Value = 72000
$ws.Range("L83").Value = 135000
$ws.Range("M83").Value = -67320
$ws.Range("N83").Value = -144360
$ws.Range("H92").Value = 373.5
$ws.Range("J92").Value = 333.33334
$ws.Range("L92").Value = 1000.00002
$ws.Range("N92").Value = -3496.00002
$ws.Range("H113").Value = 1597.8
$ws.Range("I113").Value = 995.8333
$ws.Range("K113").Value = 2987.4999
$ws.Range("M113").Value = -817.4998999999998
$ws.Range("H129").Value = 1604
$ws.Range("J129").Value = 2172.5715
$ws.Range("L129").Value = 6517.7145
$ws.Range("N129").Value = -16517.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 8999
$ws.Range("I17").Value = 8999
$ws.Range("K17").Value = 8999
$ws.Range("M17").Value = -8829
$ws.Range("H64").Value = 66444
$ws.Range("J64").Value = 66444
$ws.Range("L64").Value = 66444
$ws.Range("N64").Value = -66894
$ws.Range("H67").Value = 66444
$ws.Range("J67").Value = 66444
$ws.Range("L67").Value = 66444
$ws.Range("N67").Value = -68004
$ws.Range("H132").Value = 3412.9333
$ws.Range("I132").Value = 2654.0908
$ws.Range("J132").Value = 5499.75
$ws.Range("K132").Value = 7962.2724
$ws.Range("L132").Value = 16499.25
$ws.Range("M132").Value = -5432.2724
$ws.Range("N132").Value = -21559.25
$ws.Range("H133").Value = 56950
$ws.Range("J133").Value = 56950
$ws.Range("L133").Value = 56950
$ws.Range("N133").Value = -62010
